# Applies the scheduled-runner profit-data refresh to the Leve profit sheets.
# Each block targets one worksheet/row; values come from the updated market-price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 119.75
$ws.Range("I6").Value = 93
$ws.Range("K6").Value = 279
$ws.Range("M6").Value = -167

$ws.Range("H28").Value = 50776.39
$ws.Range("I28").Value = 1016.0769
$ws.Range("K28").Value = 1016.0769
$ws.Range("M28").Value = -531.0769

$ws.Range("H69").Value = 5503.25
$ws.Range("I69").Value = 3506.5
$ws.Range("J69").Value = 7500
$ws.Range("K69").Value = 10519.5
$ws.Range("L69").Value = 22500
$ws.Range("M69").Value = -9645.5
$ws.Range("N69").Value = -24248

$ws.Range("H72").Value = 5503.25
$ws.Range("I72").Value = 3506.5
$ws.Range("J72").Value = 7500
$ws.Range("K72").Value = 31558.5
$ws.Range("L72").Value = 67500
$ws.Range("M72").Value = -27190.5
$ws.Range("N72").Value = -76236

$ws.Range("H92").Value = 375.21738
$ws.Range("I92").Value = 311.17648
$ws.Range("J92").Value = 556.6667
$ws.Range("K92").Value = 311.17648
$ws.Range("L92").Value = 556.6667
$ws.Range("M92").Value = 936.8235199999999
$ws.Range("N92").Value = -3052.6667

$ws.Range("H103").Value = 610.2143
$ws.Range("I103").Value = 663.7143
$ws.Range("J103").Value = 556.7143
$ws.Range("K103").Value = 1991.1429
$ws.Range("L103").Value = 1670.1429
$ws.Range("M103").Value = -1405.1429
$ws.Range("N103").Value = -2842.1429

$ws.Range("H107").Value = 72065.14
$ws.Range("I107").Value = 100411
$ws.Range("K107").Value = 100411
$ws.Range("M107").Value = -98491

$ws.Range("H125").Value = 860.6842
$ws.Range("I125").Value = 697.2
$ws.Range("J125").Value = 1042.3334
$ws.Range("K125").Value = 6274.8
$ws.Range("L125").Value = 9381.000599999999
$ws.Range("M125").Value = -3814.8
$ws.Range("N125").Value = -14301.0006

$ws.Range("H135").Value = 834.5
$ws.Range("I135").Value = 546.0526
$ws.Range("K135").Value = 4914.4734
$ws.Range("M135").Value = -2379.4734

$ws.Range("H136").Value = 72869.8
$ws.Range("J136").Value = 82087.25
$ws.Range("L136").Value = 82087.25
$ws.Range("N136").Value = -92287.25

$ws.Range("H137").Value = 531132.3
$ws.Range("I137").Value = 1880
$ws.Range("K137").Value = 5640
$ws.Range("M137").Value = -3090

$ws.Range("H138").Value = 2687.1128
$ws.Range("I138").Value = 1941.36
$ws.Range("J138").Value = 3191
$ws.Range("K138").Value = 5824.08
$ws.Range("L138").Value = 9573
$ws.Range("M138").Value = -684.0799999999999
$ws.Range("N138").Value = -19853

$ws.Range("H139").Value = 99990
$ws.Range("J139").Value = 99990
$ws.Range("L139").Value = 99990
$ws.Range("N139").Value = -110270

$ws.Range("H140").Value = 91992.86
$ws.Range("J140").Value = 91992.86
$ws.Range("L140").Value = 91992.86
$ws.Range("N140").Value = -102352.86

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 128398.375
$ws.Range("I61").Value = 1491.3334
$ws.Range("K61").Value = 1491.3334
$ws.Range("M61").Value = -1279.3334

$ws.Range("H110").Value = 1376.5294
$ws.Range("I110").Value = 1136.6666
$ws.Range("J110").Value = 1952.2
$ws.Range("K110").Value = 1136.6666
$ws.Range("L110").Value = 1952.2
$ws.Range("M110").Value = 908.3334
$ws.Range("N110").Value = -6042.2

$ws.Range("H115").Value = 66200
$ws.Range("J115").Value = 80000
$ws.Range("L115").Value = 80000
$ws.Range("N115").Value = -83134

$ws.Range("H136").Value = 128398.375
$ws.Range("I136").Value = 1491.3334
$ws.Range("K136").Value = 4474.0002
$ws.Range("M136").Value = -1924.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2248.7
$ws.Range("I134").Value = 1882.7084
$ws.Range("K134").Value = 5648.1252
$ws.Range("M134").Value = -3113.1252

$ws.Range("H135").Value = 106491.664
$ws.Range("J135").Value = 106491.664
$ws.Range("L135").Value = 106491.664
$ws.Range("N135").Value = -116631.664

$ws.Range("H140").Value = 58304.082
$ws.Range("J140").Value = 58304.082
$ws.Range("L140").Value = 58304.082
$ws.Range("N140").Value = -68664.08199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 949.2083
$ws.Range("I107").Value = 957.6
$ws.Range("K107").Value = 957.6
$ws.Range("M107").Value = 962.4

$ws.Range("H127").Value = 50000
$ws.Range("J127").Value = 50000
$ws.Range("L127").Value = 50000
$ws.Range("N127").Value = -59920

$ws.Range("H138").Value = 54140
$ws.Range("J138").Value = 53657.5
$ws.Range("L138").Value = 53657.5
$ws.Range("N138").Value = -63937.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2631889
$ws.Range("I4").Value = 2830474.5
$ws.Range("K4").Value = 8491423.5
$ws.Range("M4").Value = -8491311.5

$ws.Range("H5").Value = 1464.909
$ws.Range("J5").Value = 2666
$ws.Range("L5").Value = 7998
$ws.Range("N5").Value = -8222

$ws.Range("H69").Value = 2166.6667
$ws.Range("I69").Value = 2166.6667
$ws.Range("K69").Value = 6500.000100000001
$ws.Range("M69").Value = -5689.000100000001

$ws.Range("H72").Value = 2166.6667
$ws.Range("I72").Value = 2166.6667
$ws.Range("K72").Value = 19500.0003
$ws.Range("M72").Value = -15444.0003

$ws.Range("H97").Value = 300
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 300
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 900
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -1892

$ws.Range("H135").Value = 1464.909
$ws.Range("J135").Value = 2666
$ws.Range("L135").Value = 23994
$ws.Range("N135").Value = -29064

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 78564.14
$ws.Range("J135").Value = 78564.14
$ws.Range("L135").Value = 78564.14
$ws.Range("N135").Value = -88704.14

$ws.Range("H140").Value = 43537.734
$ws.Range("J140").Value = 52437.668
$ws.Range("L140").Value = 52437.668
$ws.Range("N140").Value = -62797.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws.Range("H64").Value = 10075
$ws.Range("J64").Value = 10075
$ws.Range("L64").Value = 10075
$ws.Range("N64").Value = -10525

$ws.Range("H67").Value = 10075
$ws.Range("J67").Value = 10075
$ws.Range("L67").Value = 10075
$ws.Range("N67").Value = -11635

$ws.Range("H68").Value = 302557.28
$ws.Range("I68").Value = 352650.16
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 352650.16
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -351901.16
$ws.Range("N68").Value = -3498

$ws.Range("H71").Value = 302557.28
$ws.Range("I71").Value = 352650.16
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 1763250.8
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -1759506.8
$ws.Range("N71").Value = -17488

$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 17000
$ws.Range("I38").Value = 17000
$ws.Range("J38").Value = 17000
$ws.Range("K38").Value = 17000
$ws.Range("L38").Value = 17000
$ws.Range("M38").Value = -16527
$ws.Range("N38").Value = -17946

$ws.Range("H59").Value = 39991
$ws.Range("J59").Value = 39991
$ws.Range("L59").Value = 39991
$ws.Range("N59").Value = -41467

$ws.Range("H61").Value = 15731287
$ws.Range("I61").Value = 18352108
$ws.Range("K61").Value = 18352108
$ws.Range("M61").Value = -18351816

$ws.Range("H75").Value = 19333
$ws.Range("I75").Value = 19333
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 19333
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -18397
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 19333
$ws.Range("I78").Value = 19333
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 57999
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -53319
$ws.Range("N78").ClearContents()

$ws.Range("H100").Value = 7143449.5
$ws.Range("I100").Value = 17857390
$ws.Range("K100").Value = 35714780
$ws.Range("M100").Value = -35714239

$ws.Range("H107").Value = 9022.066000000001
$ws.Range("I107").Value = 19071
$ws.Range("K107").Value = 57213
$ws.Range("M107").Value = -55293

$ws.Range("H118").Value = 67196
$ws.Range("J118").Value = 67196
$ws.Range("L118").Value = 67196
$ws.Range("N118").Value = -70510

$ws.Range("H129").Value = 55000
$ws.Range("J129").Value = 55000
$ws.Range("L129").Value = 55000
$ws.Range("N129").Value = -65000
